$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.064.55'
$ws.Range('E2').Value = '  +3.35%  '

$ws.Range('D3').Value = '3.452.11'
$ws.Range('E3').Value = '  +2.76%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.31%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.39%  '

$ws.Range('D8').Value = '3.445.92'
$ws.Range('E8').Value = '  +2.88%  '

$ws.Range('E9').Value = '  +0.04%  '

$ws.Range('E10').Value = '  -0.63%  '

$ws.Range('E11').Value = '  +1.53%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.20'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.66%  '

$ws.Range('E13').Value = '  -1.46%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.00%  '

$ws.Range('D15').Value = '3.999.60'
$ws.Range('E15').Value = '  +2.89%  '

$ws.Range('E16').Value = '  +2.70%  '

$ws.Range('D17').Value = '3.448.91'
$ws.Range('E17').Value = '  +3.12%  '

$ws.Range('D18').Value = '67.145.62'
$ws.Range('E18').Value = '  +3.71%  '

$ws.Range('E19').Value = '  +2.49%  '

$ws.Range('E20').Value = '  -2.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.34%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '487.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.13%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +23.43%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.32%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.41%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.94'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.09%  '

$ws.Range('E29').Value = '  +4.76%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.00%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.77%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '598.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.82%  '

$ws.Range('E35').Value = '  +4.15%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.151'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.73%  '

$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.20%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.58%  '

$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.386'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.27%  '

$ws.Range('D41').Value = '3.259.69'
$ws.Range('E41').Value = '  +5.94%  '

$ws.Range('D42').Value = '0.0₃0753'
$ws.Range('E42').Value = '  +1.68%  '

$ws.Range('E43').Value = '  +5.72%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0431'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.63%  '

$ws.Range('E45').Value = '  +3.14%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +22.80%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.55%  '

$ws.Range('E48').Value = '  +0.41%  '

$ws.Range('E49').Value = '  +13.86%  '

$ws.Range('E50').Value = '  +5.48%  '

$ws.Range('E51').Value = '  +0.10%  '
